$d = $word.ActiveDocument
$startIdx = $d.Paragraphs.Count
$anchor = $d.Paragraphs($startIdx)
$r0 = $anchor.Range
$r0.Collapse(0) | Out-Null
$N = 13
for ($i = 0; $i -lt $N; $i++) {
    $r0.InsertParagraphAfter()
    $r0.Collapse(0) | Out-Null
}

# --- Paragraph 0: style=Heading1 (needs_style_set=True) ---
$idx1 = $startIdx + 1
$p1 = $d.Paragraphs($idx1)
$p1.Style = "Heading1"
$rp1 = $p1.Range
$rp1.Text = "Knärot – ekologi samt krav på livsmiljön"

# --- Paragraph 1: style=Normal (needs_style_set=True) ---
$idx2 = $startIdx + 2
$p2 = $d.Paragraphs($idx2)
$p2.Style = "Normal"
$rp2 = $p2.Range
$rp2.Text = "Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021)."

# --- Paragraph 2: style=Normal (needs_style_set=False) ---
$idx3 = $startIdx + 3
$p3 = $d.Paragraphs($idx3)
$rp3 = $p3.Range
$rp3.Text = "Samuel Johnsons doktorsavhandling "
$pEnd3 = $d.Paragraphs($idx3).Range.End - 1
$ip3 = $d.Range($pEnd3, $pEnd3)
$ip3.InsertAfter("“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“")
$ip3.Font.Italic = 1
$pEnd3 = $d.Paragraphs($idx3).Range.End - 1
$ip3 = $d.Range($pEnd3, $pEnd3)
$ip3.InsertAfter(" (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ")
$pEnd3 = $d.Paragraphs($idx3).Range.End - 1
$ip3 = $d.Range($pEnd3, $pEnd3)
$ip3.InsertAfter("“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ")
$ip3.Font.Italic = 1
$pEnd3 = $d.Paragraphs($idx3).Range.End - 1
$ip3 = $d.Range($pEnd3, $pEnd3)
$ip3.InsertAfter("Vidare ")
$pEnd3 = $d.Paragraphs($idx3).Range.End - 1
$ip3 = $d.Range($pEnd3, $pEnd3)
$ip3.InsertAfter("“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”")
$ip3.Font.Italic = 1

# --- Paragraph 3: style=Normal (needs_style_set=False) ---
$idx4 = $startIdx + 4
$p4 = $d.Paragraphs($idx4)
$rp4 = $p4.Range
$rp4.Text = "Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: "
$pEnd4 = $d.Paragraphs($idx4).Range.End - 1
$ip4 = $d.Range($pEnd4, $pEnd4)
$ip4.InsertAfter("“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”")
$ip4.Font.Italic = 1

# --- Paragraph 4: style=Normal (needs_style_set=False) ---
$idx5 = $startIdx + 5
$p5 = $d.Paragraphs($idx5)
$rp5 = $p5.Range
$rp5.Text = "En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022)."

# --- Paragraph 5: style=Normal (needs_style_set=False) ---
$idx6 = $startIdx + 6
$p6 = $d.Paragraphs($idx6)
$rp6 = $p6.Range
$rp6.Text = "Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022)."

# --- Paragraph 6: style=Heading2 (needs_style_set=True) ---
$idx7 = $startIdx + 7
$p7 = $d.Paragraphs($idx7)
$p7.Style = "Heading2"
$rp7 = $p7.Range
$rp7.Text = "Referenser - knärot"

# --- Paragraph 7: style=Normal (needs_style_set=True) ---
$idx8 = $startIdx + 8
$p8 = $d.Paragraphs($idx8)
$p8.Style = "Normal"
$rp8 = $p8.Range
$rp8.Text = "de Graaf M & Roberts M.R., 2009. "
$pEnd8 = $d.Paragraphs($idx8).Range.End - 1
$ip8 = $d.Range($pEnd8, $pEnd8)
$ip8.InsertAfter("Short-term response of the herbaceous layer within leave patches after harvest. ")
$ip8.Font.Italic = 1
$pEnd8 = $d.Paragraphs($idx8).Range.End - 1
$ip8 = $d.Range($pEnd8, $pEnd8)
$ip8.InsertAfter("Forest Ecology and Management 257, 1014-1025")

# --- Paragraph 8: style=Normal (needs_style_set=False) ---
$idx9 = $startIdx + 9
$p9 = $d.Paragraphs($idx9)
$rp9 = $p9.Range
$rp9.Text = "Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. "
$pEnd9 = $d.Paragraphs($idx9).Range.End - 1
$ip9 = $d.Range($pEnd9, $pEnd9)
$ip9.InsertAfter("Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ")
$ip9.Font.Italic = 1
$pEnd9 = $d.Paragraphs($idx9).Range.End - 1
$ip9 = $d.Range($pEnd9, $pEnd9)
$ip9.InsertAfter("Ecological Applications, 22, 2049-2064 ")

# --- Paragraph 9: style=Normal (needs_style_set=False) ---
$idx10 = $startIdx + 10
$p10 = $d.Paragraphs($idx10)
$rp10 = $p10.Range
$rp10.Text = "Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. "
$pEnd10 = $d.Paragraphs($idx10).Range.End - 1
$ip10 = $d.Range($pEnd10, $pEnd10)
$ip10.InsertAfter("Interactive effects of drought and edge exposure on old-growth forest understory species. ")
$ip10.Font.Italic = 1
$pEnd10 = $d.Paragraphs($idx10).Range.End - 1
$ip10 = $d.Range($pEnd10, $pEnd10)
$ip10.InsertAfter("Landscape Ecology, 37, sid 1839-1853")

# --- Paragraph 10: style=Normal (needs_style_set=False) ---
$idx11 = $startIdx + 11
$p11 = $d.Paragraphs($idx11)
$rp11 = $p11.Range
$rp11.Text = "Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. "
$pEnd11 = $d.Paragraphs($idx11).Range.End - 1
$ip11 = $d.Range($pEnd11, $pEnd11)
$ip11.InsertAfter("Biological legacies buffer local species extinction after logging. ")
$ip11.Font.Italic = 1
$pEnd11 = $d.Paragraphs($idx11).Range.End - 1
$ip11 = $d.Range($pEnd11, $pEnd11)
$ip11.InsertAfter("Journal of Applied Ecology. 51, 53-62.")

# --- Paragraph 11: style=Normal (needs_style_set=False) ---
$idx12 = $startIdx + 12
$p12 = $d.Paragraphs($idx12)
$rp12 = $p12.Range
$rp12.Text = "Skogsstyrelsen, 2022. "
$pEnd12 = $d.Paragraphs($idx12).Range.End - 1
$ip12 = $d.Range($pEnd12, $pEnd12)
$ip12.InsertAfter("Vägledning för hänsyn till knärot. ")
$ip12.Font.Italic = 1
$pEnd12 = $d.Paragraphs($idx12).Range.End - 1
$ip12 = $d.Range($pEnd12, $pEnd12)
$ip12.InsertAfter("https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/")

# --- Paragraph 12: style=Normal (needs_style_set=False) ---
$idx13 = $startIdx + 13
$p13 = $d.Paragraphs($idx13)
$rp13 = $p13.Range
$rp13.Text = "SLU Artdatabanken, 2021. "
$pEnd13 = $d.Paragraphs($idx13).Range.End - 1
$ip13 = $d.Range($pEnd13, $pEnd13)
$ip13.InsertAfter("Artfaktablad. Naturvård – artfakta. ")
$ip13.Font.Italic = 1
$pEnd13 = $d.Paragraphs($idx13).Range.End - 1
$ip13 = $d.Range($pEnd13, $pEnd13)
$ip13.InsertAfter("SLU Artdatabanken, Uppsala ")

# Update header date
$sec = $d.Sections(1)
$hdr = $sec.Headers(2)
$hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null